$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-type columns (B, C, E): safe to assign directly, Excel will not
# re-interpret these strings as numbers because of letters/percent/padding. ---
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +1.18%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +9.38%  '
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  +10.10%  '
$ws.Range('E13').Value = '  +8.25%  '
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('E21').Value = '  +3.94%  '
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('E24').Value = '  +23.35%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('E34').Value = '  +2.92%  '
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('E37').Value = '  -4.53%  '
$ws.Range('E38').Value = '  -3.95%  '
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('E41').Value = '  +12.49%  '
$ws.Range('E42').Value = '  +11.64%  '
$ws.Range('E43').Value = '  -4.38%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('B45').Value = 'BinanceUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  +4.79%  '
$ws.Range('B47').Value = 'MultiversX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('E47').Value = '  +10.13%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E48').Value = '  +6.14%  '
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('E51').Value = '  +1.31%  '

# --- Column D holds numeric-looking strings (thousand-dot separators, fixed
# decimals like "1.00" / "0.680") that Excel would silently coerce to real
# numbers (dropping the formatting, merging "44.099.79" into a date, etc.) if
# assigned while the cell is still General-formatted. Force the cell to Text
# format first, assign the literal string, then restore the Normal style so
# no stray 's' (style) attribute is left behind on the cell. ---
$dCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D21', 'D22', 'D23', 'D24', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D51')
$dValues = @('44.099.79', '2.354.99', '0.680', '239.21', '74.10', '0.598', '0.101', '57.31', '7.29', '0.108', '2.704.19', '16.58', '0.900', '2.355.61', '43.945.24', '6.73', '76.92', '259.37', '1.97', '2.49', '10.73', '2.27', '22.69', '176.02', '0.0761', '5.23', '5.52', '3.76', '2.35', '6.30', '0.205', '18.95', '8.95', '1.00', '4.68', '58.22', '2.51', '100.19')
for ($i = 0; $i -lt $dCells.Length; $i++) {
    $cell = $ws.Range($dCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$i]
    $cell.Style = "Normal"
}
